# Update the timestamps recorded in the "handback-status" report.
# The workbook has three sheets: "Overview", "zh-cn", "de-de".
# Several cells hold datetime strings (stored as text, not Excel date
# serials) recording when xliff files were generated / handed back.
# This script bumps those timestamps forward, as described by the diff.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G2 ("Latest HO Xliff Generate Date") and de-de!H2 ("Correspond
# Handoff Datetime") share the same original value 2016-08-17 13:02:29.
$wsOverview.Range("G2").Value = "2016-08-17 13:03:20"
$wsDeDe.Range("H2").Value = "2016-08-17 13:03:20"

# zh-cn!H2 ("Correspond Handoff Datetime") and zh-cn!K2 ("Correspond
# Handback DateTime").
$wsZhCn.Range("H2").Value = "2016-08-17 13:03:15"
$wsZhCn.Range("K2").Value = "2016-08-17 13:03:32"

# de-de!K2 ("Correspond Handback DateTime").
$wsDeDe.Range("K2").Value = "2016-08-17 13:03:39"
